$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32 (ALC)
$ws.Range("H32").Value = 3124.9167
$ws.Range("I32").Value = 1777.6875
$ws.Range("J32").Value = 4202.7
$ws.Range("K32").Value = 1777.6875
$ws.Range("L32").Value = 4202.7
$ws.Range("M32").Value = -1451.6875
$ws.Range("N32").Value = -4854.7

# Row 107 (ALC)
$ws.Range("H107").Value = 258.0909
$ws.Range("I107").Value = 258.0909
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 258.0909
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1661.9091
$ws.Range("N107").ClearContents()

# Row 141 (ALC)
$ws.Range("H141").Value = 7206.4585
$ws.Range("I141").Value = 2607.2222
$ws.Range("J141").Value = 21004.166
$ws.Range("K141").Value = 7821.6666
$ws.Range("L141").Value = 63012.49800000001
$ws.Range("M141").Value = -2641.6666
$ws.Range("N141").Value = -73372.49800000001

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (ARM)
$ws.Range("H32").Value = 11243.427
$ws.Range("I32").Value = 3270.1272
$ws.Range("K32").Value = 3270.1272
$ws.Range("M32").Value = -2983.1272

# Row 45 (ARM)
$ws.Range("H45").Value = 2021.3793
$ws.Range("I45").Value = 2015.619
$ws.Range("J45").Value = 2036.5
$ws.Range("K45").Value = 2015.619
$ws.Range("L45").Value = 2036.5
$ws.Range("M45").Value = -1638.619
$ws.Range("N45").Value = -2790.5

# Row 61 (ARM)
$ws.Range("H61").Value = 1347.7778
$ws.Range("I61").Value = 1002.5
$ws.Range("J61").Value = 2038.3334
$ws.Range("K61").Value = 1002.5
$ws.Range("L61").Value = 2038.3334
$ws.Range("M61").Value = -790.5
$ws.Range("N61").Value = -2462.3334

# Row 105 (ARM)
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

# Row 122 (ARM)
$ws.Range("H122").Value = 2090.4644
$ws.Range("I122").Value = 1978.7727
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 5936.3181
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -3486.3181
$ws.Range("N122").Value = -12400

# Row 132 (ARM)
$ws.Range("H132").Value = 1912.3226
$ws.Range("I132").Value = 1448.1666
$ws.Range("J132").Value = 3503.7144
$ws.Range("K132").Value = 4344.4998
$ws.Range("L132").Value = 10511.1432
$ws.Range("M132").Value = -1814.4998
$ws.Range("N132").Value = -15571.1432

# Row 136 (ARM)
$ws.Range("H136").Value = 1347.7778
$ws.Range("I136").Value = 1002.5
$ws.Range("J136").Value = 2038.3334
$ws.Range("K136").Value = 3007.5
$ws.Range("L136").Value = 6115.0002
$ws.Range("M136").Value = -457.5
$ws.Range("N136").Value = -11215.0002

$ws = $wb.Worksheets.Item("BSM")
# Row 20 (BSM)
$ws.Range("H20").Value = 2081.16
$ws.Range("I20").Value = 1744.9286
$ws.Range("J20").Value = 2509.0908
$ws.Range("K20").Value = 1744.9286
$ws.Range("L20").Value = 2509.0908
$ws.Range("M20").Value = -1497.9286
$ws.Range("N20").Value = -3003.0908

# Row 63 (BSM)
$ws.Range("H63").Value = 42000
$ws.Range("J63").Value = 42000
$ws.Range("L63").Value = 42000
$ws.Range("N63").Value = -43372

# Row 66 (BSM)
$ws.Range("H66").Value = 42000
$ws.Range("J66").Value = 42000
$ws.Range("L66").Value = 126000
$ws.Range("N66").Value = -132864

# Row 86 (BSM)
$ws.Range("H86").Value = 25003202
$ws.Range("I86").Value = 33335850
$ws.Range("J86").Value = 5253.5
$ws.Range("K86").Value = 33335850
$ws.Range("L86").Value = 5253.5
$ws.Range("M86").Value = -33334727
$ws.Range("N86").Value = -7499.5

# Row 89 (BSM)
$ws.Range("H89").Value = 25003202
$ws.Range("I89").Value = 33335850
$ws.Range("J89").Value = 5253.5
$ws.Range("K89").Value = 166679250
$ws.Range("L89").Value = 26267.5
$ws.Range("M89").Value = -166673634
$ws.Range("N89").Value = -37499.5

# Row 94 (BSM)
$ws.Range("H94").Value = 12859.8125
$ws.Range("I94").Value = 362.69232
$ws.Range("J94").Value = 67014
$ws.Range("K94").Value = 362.69232
$ws.Range("L94").Value = 67014
$ws.Range("M94").Value = 88.30768
$ws.Range("N94").Value = -67916

# Row 134 (BSM)
$ws.Range("H134").Value = 2242.6924
$ws.Range("I134").Value = 1535.5
$ws.Range("J134").Value = 4600
$ws.Range("K134").Value = 4606.5
$ws.Range("L134").Value = 13800
$ws.Range("M134").Value = -2071.5
$ws.Range("N134").Value = -18870

$ws = $wb.Worksheets.Item("CRP")
# Row 106 (CRP)
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 17 (CUL)
$ws.Range("H17").Value = 300
$ws.Range("I17").Value = 300
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 900
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -731
$ws.Range("N17").ClearContents()

# Row 59 (CUL)
$ws.Range("H59").Value = 2625
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 2625
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 7875
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -8955

$ws = $wb.Worksheets.Item("GSM")
# Row 24 (GSM)
$ws.Range("H24").Value = 10406005
$ws.Range("I24").Value = 26000004
$ws.Range("K24").Value = 26000004
$ws.Range("M24").Value = -25999831

$ws = $wb.Worksheets.Item("LTW")
# Row 16 (LTW)
$ws.Range("H16").Value = 1286.3158
$ws.Range("I16").Value = 775.93335
$ws.Range("J16").Value = 3200.25
$ws.Range("K16").Value = 775.93335
$ws.Range("L16").Value = 3200.25
$ws.Range("M16").Value = -605.93335
$ws.Range("N16").Value = -3540.25

# Row 22 (LTW)
$ws.Range("H22").Value = 1588025.4
$ws.Range("I22").Value = 2381431
$ws.Range("J22").Value = 1214.2858
$ws.Range("K22").Value = 2381431
$ws.Range("L22").Value = 1214.2858
$ws.Range("M22").Value = -2381136
$ws.Range("N22").Value = -1804.2858

# Row 27 (LTW)
$ws.Range("H27").Value = 1588025.4
$ws.Range("I27").Value = 2381431
$ws.Range("J27").Value = 1214.2858
$ws.Range("K27").Value = 2381431
$ws.Range("L27").Value = 1214.2858
$ws.Range("M27").Value = -2381324
$ws.Range("N27").Value = -1428.2858

# Row 40 (LTW)
$ws.Range("H40").Value = 1258.4615
$ws.Range("I40").Value = 997.5
$ws.Range("J40").Value = 2128.3333
$ws.Range("K40").Value = 997.5
$ws.Range("L40").Value = 2128.3333
$ws.Range("M40").Value = -861.5
$ws.Range("N40").Value = -2400.3333

# Row 93 (LTW)
$ws.Range("H93").Value = 972.5454999999999
$ws.Range("I93").Value = 1078.0714
$ws.Range("J93").Value = 787.875
$ws.Range("K93").Value = 1078.0714
$ws.Range("L93").Value = 787.875
$ws.Range("M93").Value = 169.9286
$ws.Range("N93").Value = -3283.875

# Row 105 (LTW)
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 132 (WVR)
$ws.Range("H132").Value = 2088.074
$ws.Range("I132").Value = 1600.7333
$ws.Range("J132").Value = 2697.25
$ws.Range("K132").Value = 4802.199900000001
$ws.Range("L132").Value = 8091.75
$ws.Range("M132").Value = -2272.199900000001
$ws.Range("N132").Value = -13151.75
